# Update the "想去人数" (interest count) figures in the F column on the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to reflect newly
# published output.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 84
$ws1.Range("F3").Value = 811
$ws1.Range("F6").Value = 119
$ws1.Range("F7").Value = 345
$ws1.Range("F8").Value = 4603
$ws1.Range("F10").Value = 5043
$ws1.Range("F11").Value = 574
$ws1.Range("F12").Value = 1269
$ws1.Range("F13").Value = 90

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 84
$ws4.Range("F3").Value = 811
$ws4.Range("F6").Value = 119
$ws4.Range("F8").Value = 345
$ws4.Range("F9").Value = 4603
$ws4.Range("F11").Value = 5043
$ws4.Range("F12").Value = 574
$ws4.Range("F13").Value = 1269
$ws4.Range("F14").Value = 90
